# Applies the numeric cell updates recorded in the commit diff for
# Sheets/Zeromus_Profits.xlsx, across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# worksheets of this workbook.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4224.3477
$ws.Range("I17").Value = 780
$ws.Range("J17").Value = 4552.381
$ws.Range("K17").Value = 2340
$ws.Range("L17").Value = 13657.143
$ws.Range("M17").Value = -2172
$ws.Range("N17").Value = -13993.143
$ws.Range("H64").Value = 2868.4546
$ws.Range("J64").Value = 2774.2307
$ws.Range("L64").Value = 2774.2307
$ws.Range("N64").Value = -3270.2307
$ws.Range("H67").Value = 2868.4546
$ws.Range("J67").Value = 2774.2307
$ws.Range("L67").Value = 2774.2307
$ws.Range("N67").Value = -4490.2307
$ws.Range("H70").Value = 1698.4839
$ws.Range("I70").Value = 1811.6538
$ws.Range("K70").Value = 5434.9614
$ws.Range("M70").Value = -5164.9614
$ws.Range("H73").Value = 1698.4839
$ws.Range("I73").Value = 1811.6538
$ws.Range("K73").Value = 5434.9614
$ws.Range("M73").Value = -4498.9614
$ws.Range("H98").Value = 518.9583
$ws.Range("I98").Value = 518.9583
$ws.Range("K98").Value = 518.9583
$ws.Range("M98").Value = 979.0417
$ws.Range("H122").Value = 518.9583
$ws.Range("I122").Value = 518.9583
$ws.Range("K122").Value = 1556.8749
$ws.Range("M122").Value = 893.1251
$ws.Range("H129").Value = 935.9643
$ws.Range("I129").Value = 321.75
$ws.Range("J129").Value = 1038.3334
$ws.Range("K129").Value = 965.25
$ws.Range("L129").Value = 3115.0002
$ws.Range("M129").Value = 4034.75
$ws.Range("N129").Value = -13115.0002
$ws.Range("H130").Value = 39668.57
$ws.Range("J130").Value = 39668.57
$ws.Range("L130").Value = 39668.57
$ws.Range("N130").Value = -49708.57
$ws.Range("H131").Value = 965
$ws.Range("I131").Value = 211.875
$ws.Range("J131").Value = 1969.1666
$ws.Range("K131").Value = 635.625
$ws.Range("L131").Value = 5907.4998
$ws.Range("M131").Value = 4404.375
$ws.Range("N131").Value = -15987.4998
$ws.Range("H132").Value = 3828.4897
$ws.Range("I132").Value = 3833.262
$ws.Range("J132").Value = 3799.8572
$ws.Range("K132").Value = 11499.786
$ws.Range("L132").Value = 11399.5716
$ws.Range("M132").Value = -8969.786
$ws.Range("N132").Value = -16459.5716
$ws.Range("H133").Value = 41593.332
$ws.Range("J133").Value = 41593.332
$ws.Range("L133").Value = 41593.332
$ws.Range("N133").Value = -51713.332
$ws.Range("H138").Value = 2989.9885
$ws.Range("I138").Value = 1445.6923
$ws.Range("J138").Value = 3261.2837
$ws.Range("K138").Value = 4337.0769
$ws.Range("L138").Value = 9783.8511
$ws.Range("M138").Value = 802.9231
$ws.Range("N138").Value = -20063.8511

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 835.3913
$ws.Range("I2").Value = 824.2632
$ws.Range("K2").Value = 824.2632
$ws.Range("M2").Value = -711.2632
$ws.Range("H32").Value = 12255.17
$ws.Range("I32").Value = 3878.8142
$ws.Range("J32").Value = 31800
$ws.Range("K32").Value = 3878.8142
$ws.Range("L32").Value = 31800
$ws.Range("M32").Value = -3591.8142
$ws.Range("N32").Value = -32374
$ws.Range("H63").Value = 2698.524
$ws.Range("J63").Value = 3845.5715
$ws.Range("L63").Value = 3845.5715
$ws.Range("N63").Value = -5217.5715
$ws.Range("H66").Value = 2698.524
$ws.Range("J66").Value = 3845.5715
$ws.Range("L66").Value = 19227.8575
$ws.Range("N66").Value = -26091.8575
$ws.Range("H74").Value = 4928.1055
$ws.Range("I74").Value = 4846.3335
$ws.Range("J74").Value = 5068.2856
$ws.Range("K74").Value = 4846.3335
$ws.Range("L74").Value = 5068.2856
$ws.Range("M74").Value = -3972.3335
$ws.Range("N74").Value = -6816.2856
$ws.Range("H77").Value = 4928.1055
$ws.Range("I77").Value = 4846.3335
$ws.Range("J77").Value = 5068.2856
$ws.Range("K77").Value = 24231.6675
$ws.Range("L77").Value = 25341.428
$ws.Range("M77").Value = -19863.6675
$ws.Range("N77").Value = -34077.428
$ws.Range("H102").Value = 1236.25
$ws.Range("J102").Value = 1500
$ws.Range("L102").Value = 1500
$ws.Range("N102").Value = -4744
$ws.Range("H116").Value = 835.3913
$ws.Range("I116").Value = 824.2632
$ws.Range("K116").Value = 824.2632
$ws.Range("M116").Value = 1469.7368
$ws.Range("H122").Value = 1427.8182
$ws.Range("I122").Value = 954.2857
$ws.Range("J122").Value = 2256.5
$ws.Range("K122").Value = 2862.8571
$ws.Range("L122").Value = 6769.5
$ws.Range("M122").Value = -412.8571000000002
$ws.Range("N122").Value = -11669.5
$ws.Range("H132").Value = 2094.0881
$ws.Range("I132").Value = 1276.8334
$ws.Range("J132").Value = 3013.5
$ws.Range("K132").Value = 3830.5002
$ws.Range("L132").Value = 9040.5
$ws.Range("M132").Value = -1300.5002
$ws.Range("N132").Value = -14100.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 835.3913
$ws.Range("I3").Value = 824.2632
$ws.Range("K3").Value = 824.2632
$ws.Range("M3").Value = -710.2632
$ws.Range("H134").Value = 1347.5
$ws.Range("I134").Value = 1050.7838
$ws.Range("J134").Value = 2915.8572
$ws.Range("K134").Value = 3152.3514
$ws.Range("L134").Value = 8747.571599999999
$ws.Range("M134").Value = -617.3513999999996
$ws.Range("N134").Value = -13817.5716

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3280789.8
$ws.Range("I31").Value = 6062112.5
$ws.Range("K31").Value = 6062112.5
$ws.Range("M31").Value = -6061817.5
$ws.Range("H34").Value = 3280789.8
$ws.Range("I34").Value = 6062112.5
$ws.Range("K34").Value = 6062112.5
$ws.Range("M34").Value = -6061910.5
$ws.Range("H42").Value = 10000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 10000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 10000
$ws.Range("M42").ClearContents()
$ws.Range("N42").Value = -11186
$ws.Range("H114").Value = 26750
$ws.Range("J114").Value = 26750
$ws.Range("L114").Value = 26750
$ws.Range("N114").Value = -35428
$ws.Range("H119").Value = 34112.2
$ws.Range("I119").Value = 16666.666
$ws.Range("J119").Value = 60280.5
$ws.Range("K119").Value = 16666.666
$ws.Range("L119").Value = 60280.5
$ws.Range("M119").Value = -11828.666
$ws.Range("N119").Value = -69956.5
$ws.Range("H122").Value = 3175574.5
$ws.Range("I122").Value = 4167553.2
$ws.Range("J122").Value = 1242.8
$ws.Range("K122").Value = 12502659.6
$ws.Range("L122").Value = 3728.4
$ws.Range("M122").Value = -12500209.6
$ws.Range("N122").Value = -8628.4
$ws.Range("H134").Value = 2870.4
$ws.Range("I134").Value = 2468.4358
$ws.Range("K134").Value = 7405.307400000001
$ws.Range("M134").Value = -4870.307400000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2755.6365
$ws.Range("I109").Value = 1303.25
$ws.Range("J109").Value = 3585.5715
$ws.Range("K109").Value = 3909.75
$ws.Range("L109").Value = 10756.7145
$ws.Range("M109").Value = -2869.75
$ws.Range("N109").Value = -12836.7145
$ws.Range("H121").Value = 1197.0541
$ws.Range("J121").Value = 1763.174
$ws.Range("L121").Value = 5289.522
$ws.Range("N121").Value = -7909.522
$ws.Range("H131").Value = 2223281
$ws.Range("I131").Value = 8333818
$ws.Range("J131").Value = 1267.5
$ws.Range("K131").Value = 25001454
$ws.Range("L131").Value = 3802.5
$ws.Range("M131").Value = -24996414
$ws.Range("N131").Value = -13882.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1477.6
$ws.Range("I102").Value = 1472
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 1472
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 150
$ws.Range("N102").Value = -4744
$ws.Range("H126").Value = 2577.3333
$ws.Range("I126").Value = 2251.7144
$ws.Range("J126").Value = 3228.5715
$ws.Range("K126").Value = 6755.1432
$ws.Range("L126").Value = 9685.7145
$ws.Range("M126").Value = -4285.1432
$ws.Range("N126").Value = -14625.7145

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2072.238
$ws.Range("I7").Value = 1875.2222
$ws.Range("J7").Value = 2220
$ws.Range("K7").Value = 1875.2222
$ws.Range("L7").Value = 2220
$ws.Range("M7").Value = -1763.2222
$ws.Range("N7").Value = -2444
$ws.Range("H126").Value = 2072.238
$ws.Range("I126").Value = 1875.2222
$ws.Range("J126").Value = 2220
$ws.Range("K126").Value = 5625.6666
$ws.Range("L126").Value = 6660
$ws.Range("M126").Value = -3155.6666
$ws.Range("N126").Value = -11600

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 644.4231
$ws.Range("I107").Value = 433.2353
$ws.Range("J107").Value = 1043.3334
$ws.Range("K107").Value = 1299.7059
$ws.Range("L107").Value = 3130.0002
$ws.Range("M107").Value = 620.2941000000001
$ws.Range("N107").Value = -6970.0002

